$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-17 (Sending cluster x Target cluster combinations)
# Columns: RowNum, A (Sending cluster), D (Target cluster), E..T (numeric metrics)
$rows = @(
    @(2, "ECs", "ECs", 3, 1, 3.588675333333333, 10.766026, 0.1217029912931362, 0.1217029912931362, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 0.3991492177248889, 3.592342959524, 0.0009745196821405287, 0.0009745196821405286),
    @(3, "ECs", "FAPs", 3, 1, 3.588675333333333, 10.766026, 0.1217029912931362, 0.1217029912931362, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 42.43690403464267, 381.932136311784, 0.1036093681119831, 0.1036093681119831),
    @(4, "ECs", "M2", 3, 1, 3.588675333333333, 10.766026, 0.1217029912931362, 0.1217029912931362, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 0.1361447723457778, 1.225302951112, 0.0003323963930776082, 0.0003323963930776082),
    @(5, "ECs", "sCs", 3, 1, 3.588675333333333, 10.766026, 0.1217029912931362, 0.1217029912931362, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 6.875593312588, 61.880339813292, 0.01678670710593503, 0.01678670710593503),
    @(6, "FAPs", "ECs", 3, 1, 18.16892433333334, 54.50677300000001, 0.6161639698655711, 0.6161639698655711, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 2.02083255266689, 18.18749297400201, 0.004933846815757826, 0.004933846815757825),
    @(7, "FAPs", "FAPs", 3, 1, 18.16892433333334, 54.50677300000001, 0.6161639698655711, 0.6161639698655711, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 214.8516727564147, 1933.665054807732, 0.5245586726572369, 0.5245586726572368),
    @(8, "FAPs", "M2", 3, 1, 18.16892433333334, 54.50677300000001, 0.6161639698655711, 0.6161639698655711, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 0.6892805387417779, 6.203524848676001, 0.001682873025153382, 0.001682873025153382),
    @(9, "FAPs", "sCs", 3, 1, 18.16892433333334, 54.50677300000001, 0.6161639698655711, 0.6161639698655711, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 34.81009649517401, 313.2908684565661, 0.08498857736742303, 0.08498857736742303),
    @(10, "M2", "ECs", 3, 1, 2.554717666666666, 7.664153, 0.08663831443731085, 0.08663831443731086, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 0.2841476209024444, 2.557328588122, 0.0006937441861496878, 0.0006937441861496878),
    @(11, "M2", "FAPs", 3, 1, 2.554717666666666, 7.664153, 0.08663831443731085, 0.08663831443731086, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 30.21011888396133, 271.891069955652, 0.07375776813501651, 0.07375776813501651),
    @(12, "M2", "M2", 3, 1, 2.554717666666666, 7.664153, 0.08663831443731085, 0.08663831443731086, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 0.09691917569288887, 0.8722725812359999, 0.0002366274067325241, 0.0002366274067325241),
    @(13, "M2", "sCs", 3, 1, 2.554717666666666, 7.664153, 0.08663831443731085, 0.08663831443731086, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 4.894619343614, 44.051574092526, 0.01195017470941212, 0.01195017470941212),
    @(14, "sCs", "ECs", 3, 1, 5.174840666666667, 15.524522, 0.1754947244039818, 0.1754947244039818, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 0.5755699282031113, 5.180129353828001, 0.001405249462041393, 0.001405249462041392),
    @(15, "sCs", "FAPs", 3, 1, 5.174840666666667, 15.524522, 0.1754947244039818, 0.1754947244039818, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 61.19367074700535, 550.743036723048, 0.1494038668177636, 0.1494038668177635),
    @(16, "sCs", "M2", 3, 1, 5.174840666666667, 15.524522, 0.1754947244039818, 0.1754947244039818, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 0.1963196553182222, 1.766876897864, 0.0004793128975402787, 0.0004793128975402787),
    @(17, "sCs", "sCs", 3, 1, 5.174840666666667, 15.524522, 0.1754947244039818, 0.1754947244039818, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 9.914549681036002, 89.23094712932401, 0.02420629522663653, 0.02420629522663653)
)

foreach ($r in $rows) {
    $rownum = $r[0]
    $ws.Cells.Item($rownum, 1).Value = $r[1]
    $ws.Cells.Item($rownum, 2).Value = "Tfpi"
    $ws.Cells.Item($rownum, 3).Value = "Vldlr"
    $ws.Cells.Item($rownum, 4).Value = $r[2]
    $ws.Cells.Item($rownum, 5).Value = $r[3]
    $ws.Cells.Item($rownum, 6).Value = $r[4]
    $ws.Cells.Item($rownum, 7).Value = $r[5]
    $ws.Cells.Item($rownum, 8).Value = $r[6]
    $ws.Cells.Item($rownum, 9).Value = $r[7]
    $ws.Cells.Item($rownum, 10).Value = $r[8]
    $ws.Cells.Item($rownum, 11).Value = $r[9]
    $ws.Cells.Item($rownum, 12).Value = $r[10]
    $ws.Cells.Item($rownum, 13).Value = $r[11]
    $ws.Cells.Item($rownum, 14).Value = $r[12]
    $ws.Cells.Item($rownum, 15).Value = $r[13]
    $ws.Cells.Item($rownum, 16).Value = $r[14]
    $ws.Cells.Item($rownum, 17).Value = $r[15]
    $ws.Cells.Item($rownum, 18).Value = $r[16]
    $ws.Cells.Item($rownum, 19).Value = $r[17]
    $ws.Cells.Item($rownum, 20).Value = $r[18]
}

Write-Host "Updated rows 2-17 with new data"
